$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C8").Value = 1324
$ws.Range("D8").Value = 210
$ws.Range("E8").Value = 1114
$ws.Range("F8").Value = 8.613617719442166
$ws.Range("G8").Value = 84.13897280966768
$ws.Range("H8").Value = 15.86102719033233
